$d = $word.ActiveDocument

# 1) Merge the split runs " " + "(característica)" into one run for the
#    "idade;" bullet (the "nome;" bullet a few paragraphs above already
#    has its runs merged and must stay untouched). Scope the Find to that
#    single paragraph so the other "(característica)" text isn't touched.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*idade;*(característica)*") {
        $para.Range.Find.Execute(
            " (característica)", $true, $false, $false, $false, $false,
            $true, 1, $false, " (característica)", 1) | Out-Null
    }
}

# 2) Bump line spacing from single (240 twips / auto) to 1.5 lines
#    (360 twips / auto) for the closing "}" and the following
#    "const eduardo = new Pessoa();" paragraph (the only two paragraphs
#    still left at single spacing).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Format.LineSpacingRule -eq 0) {
        $para.Format.LineSpacingRule = 1
    }
}

# 3) Merge " = " + "nome" + ";" runs into a single run (leave the
#    preceding "this.nome" run — wrapped in its own proofErr markers —
#    untouched).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*this.nome = nome;*") {
        $para.Range.Find.Execute(
            " = nome;", $true, $false, $false, $false, $false,
            $true, 1, $false, " = nome;", 1) | Out-Null
    }
}

# 4) Merge " = " + "idade" + ";" runs into a single run (leave the
#    preceding "this.idade" run untouched).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*this.idade = idade;*") {
        $para.Range.Find.Execute(
            " = idade;", $true, $false, $false, $false, $false,
            $true, 1, $false, " = idade;", 1) | Out-Null
    }
}

# 5) Merge " + clicar no " + "ínicio" + " de cada linha que se deseja
#    selecionar" runs into a single run.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*clicar no*") {
        $para.Range.Find.Execute(
            " + clicar no ínicio de cada linha que se deseja selecionar",
            $true, $false, $false, $false, $false, $true, 1, $false,
            " + clicar no ínicio de cada linha que se deseja selecionar",
            1) | Out-Null
    }
}
